$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.351.41"
$ws.Range("E2").Value = "  -1.69%  "

$ws.Range("D3").Value = "2.536.47"
$ws.Range("E3").Value = "  -0.80%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.87"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.97"
$ws.Range("E6").Value = "  +2.29%  "

$ws.Range("E7").Value = "  -1.17%  "

$ws.Range("E8").Value = "  +0.14%  "

$ws.Range("E9").Value = "  -3.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.55"
$ws.Range("E10").Value = "  -0.19%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0802"
$ws.Range("E11").Value = "  -1.19%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.34"
$ws.Range("E12").Value = "  -1.78%  "

$ws.Range("E13").Value = "  -0.15%  "

$ws.Range("D14").Value = "2.927.37"
$ws.Range("E14").Value = "  -0.73%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.93"
$ws.Range("E15").Value = "  +6.05%  "

$ws.Range("D16").Value = "2.596.15"
$ws.Range("E16").Value = "  +3.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.835"
$ws.Range("E17").Value = "  -1.11%  "

$ws.Range("D18").Value = "42.367.99"
$ws.Range("E18").Value = "  -1.72%  "

$ws.Range("E19").Value = "  -1.27%  "

$ws.Range("D20").Value = "0.0₃0946"
$ws.Range("E20").Value = "  -1.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.16"
$ws.Range("E21").Value = "  -3.39%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.75"
$ws.Range("E22").Value = "  -1.68%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.89"
$ws.Range("E23").Value = "  -4.52%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.89"
$ws.Range("E24").Value = "  -1.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.04"
$ws.Range("E25").Value = "  -1.08%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.31"
$ws.Range("E27").Value = "  -1.37%  "

$ws.Range("E28").Value = "  -4.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.98"
$ws.Range("E29").Value = "  -0.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.11"
$ws.Range("E30").Value = "  -1.41%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.89"
$ws.Range("E31").Value = "  +1.85%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.71"
$ws.Range("E32").Value = "  -1.97%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.84"
$ws.Range("E33").Value = "  +16.86%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0797"
$ws.Range("E34").Value = "  -0.71%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.62"
$ws.Range("E35").Value = "  -3.17%  "

$ws.Range("E36").Value = "  -4.27%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.15"
$ws.Range("E37").Value = "  -5.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.95"
$ws.Range("E38").Value = "  -6.65%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.111"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.117"
$ws.Range("E40").Value = "  -0.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.22"
$ws.Range("E41").Value = "  +9.42%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.73"
$ws.Range("E42").Value = "  -2.63%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.28"
$ws.Range("E44").Value = "  +1.33%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0295"
$ws.Range("E45").Value = "  -2.74%  "

$ws.Range("D46").Value = "1.952.25"
$ws.Range("E46").Value = "  -1.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.90"
$ws.Range("E47").Value = "  -1.28%  "

$ws.Range("D48").Value = "2.783.56"
$ws.Range("E48").Value = "  -0.59%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "80.71"
$ws.Range("E49").Value = "  -4.55%  "

$ws.Range("E50").Value = "  -1.50%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.03"
$ws.Range("E51").Value = "  -2.88%  "

Write-Host "Updated cryptos list"
